# Weekly update: add a new week's worth of "Angeleno" price records
# (Comercializadora del Agro de Limarí - Ciruela) at the top of the
# data block (rows 12-14), pushing all existing records down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 12 (existing rows 12.. shift down to 15..)
$ws.Range("A12:T14").EntireRow.Insert()

# Row 12: Angeleno / Especial
$ws.Range("A12").Value = 2
$ws.Range("B12").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 44623
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100103
$ws.Range("H12").Value = "Frutos de hueso (carozo)"
$ws.Range("I12").Value = 100103002
$ws.Range("J12").Value = "Ciruela"
$ws.Range("K12").Value = "Angeleno"
$ws.Range("L12").Value = "Especial"
$ws.Range("M12").Value = 16
$ws.Range("N12").Value = 235000
$ws.Range("O12").Value = 240000
$ws.Range("P12").Value = 237500
$ws.Range("Q12").Value = "$/bins (450 kilos)"
$ws.Range("R12").Value = "Región Metropolitana"
$ws.Range("S12").Value = 528
$ws.Range("T12").Value = 450

# Row 13: Angeleno / Primera
$ws.Range("A13").Value = 2
$ws.Range("B13").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44623
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100103
$ws.Range("H13").Value = "Frutos de hueso (carozo)"
$ws.Range("I13").Value = 100103002
$ws.Range("J13").Value = "Ciruela"
$ws.Range("K13").Value = "Angeleno"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 16
$ws.Range("N13").Value = 195000
$ws.Range("O13").Value = 200000
$ws.Range("P13").Value = 197500
$ws.Range("Q13").Value = "$/bins (450 kilos)"
$ws.Range("R13").Value = "Región Metropolitana"
$ws.Range("S13").Value = 439
$ws.Range("T13").Value = 450

# Row 14: Angeleno / Segunda
$ws.Range("A14").Value = 2
$ws.Range("B14").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44623
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100103
$ws.Range("H14").Value = "Frutos de hueso (carozo)"
$ws.Range("I14").Value = 100103002
$ws.Range("J14").Value = "Ciruela"
$ws.Range("K14").Value = "Angeleno"
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 16
$ws.Range("N14").Value = 155000
$ws.Range("O14").Value = 160000
$ws.Range("P14").Value = 157500
$ws.Range("Q14").Value = "$/bins (450 kilos)"
$ws.Range("R14").Value = "Región Metropolitana"
$ws.Range("S14").Value = 350
$ws.Range("T14").Value = 450

# Make sure the date cells use the same date style as the rest of column D
$ws.Range("D12:D14").NumberFormat = $ws.Range("D15").NumberFormat
